# ----------------------------------------------------------------------------
# AMT_V3 Test Report - "Functionality Test" sheet
# Fills in the previously-blank magnet/coil voltage readings (columns C:F) for
# the 5 functionality-test samples (rows 19-24, 26-31, 33-38, 40-45, 47-52).
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functionality Test")

# Sample 1 (rows 19-24)
$ws.Range("C19").Value = 3.5169999999999999
$ws.Range("D19").Value = 1.992
$ws.Range("E19").Value = 1.992
$ws.Range("F19").Value = 3.512
$ws.Range("C20").Value = 3.516
$ws.Range("D20").Value = 3.5150000000000001
$ws.Range("E20").Value = 1.9910000000000001
$ws.Range("F20").Value = 1.992
$ws.Range("C21").Value = 1.998
$ws.Range("D21").Value = 3.5150000000000001
$ws.Range("E21").Value = 3.5110000000000001
$ws.Range("F21").Value = 1.992
$ws.Range("C22").Value = 1.9970000000000001
$ws.Range("D22").Value = 1.992
$ws.Range("E22").Value = 3.5110000000000001
$ws.Range("F22").Value = 3.5110000000000001
$ws.Range("C23").Value = 3.516
$ws.Range("D23").Value = 3.5139999999999998
$ws.Range("E23").Value = 3.512
$ws.Range("F23").Value = 3.512
$ws.Range("C24").Value = 1.9970000000000001
$ws.Range("D24").Value = 3.516
$ws.Range("E24").Value = 1.9910000000000001
$ws.Range("F24").Value = 3.512

# Sample 2 (rows 26-31)
$ws.Range("C26").Value = 3.492
$ws.Range("D26").Value = 1.986
$ws.Range("E26").Value = 1.9950000000000001
$ws.Range("F26").Value = 3.4910000000000001
$ws.Range("C27").Value = 3.4929999999999999
$ws.Range("D27").Value = 3.49
$ws.Range("E27").Value = 1.996
$ws.Range("F27").Value = 1.9890000000000001
$ws.Range("C28").Value = 1.99
$ws.Range("D28").Value = 3.4910000000000001
$ws.Range("E28").Value = 3.496
$ws.Range("F28").Value = 1.99
$ws.Range("C29").Value = 1.99
$ws.Range("D29").Value = 1.986
$ws.Range("E29").Value = 3.4969999999999999
$ws.Range("F29").Value = 3.4940000000000002
$ws.Range("C30").Value = 3.4929999999999999
$ws.Range("D30").Value = 3.4910000000000001
$ws.Range("E30").Value = 3.4969999999999999
$ws.Range("F30").Value = 3.4929999999999999
$ws.Range("C31").Value = 1.9910000000000001
$ws.Range("D31").Value = 3.4929999999999999
$ws.Range("E31").Value = 1.996
$ws.Range("F31").Value = 3.4950000000000001

# Sample 3 (rows 33-38)
$ws.Range("C33").Value = 3.5089999999999999
$ws.Range("D33").Value = 1.9890000000000001
$ws.Range("E33").Value = 1.9990000000000001
$ws.Range("F33").Value = 3.512
$ws.Range("C34").Value = 3.508
$ws.Range("D34").Value = 3.5059999999999998
$ws.Range("E34").Value = 1.998
$ws.Range("F34").Value = 1.9970000000000001
$ws.Range("C35").Value = 1.996
$ws.Range("D35").Value = 3.5049999999999999
$ws.Range("E35").Value = 3.51
$ws.Range("F35").Value = 1.998
$ws.Range("C36").Value = 1.996
$ws.Range("D36").Value = 1.9890000000000001
$ws.Range("E36").Value = 3.51
$ws.Range("F36").Value = 3.5110000000000001
$ws.Range("C37").Value = 3.5089999999999999
$ws.Range("D37").Value = 3.5059999999999998
$ws.Range("E37").Value = 3.51
$ws.Range("F37").Value = 3.512
$ws.Range("C38").Value = 1.996
$ws.Range("D38").Value = 3.5070000000000001
$ws.Range("E38").Value = 2
$ws.Range("F38").Value = 3.5139999999999998

# Sample 4 (rows 40-45)
$ws.Range("C40").Value = 3.5070000000000001
$ws.Range("D40").Value = 1.99
$ws.Range("E40").Value = 1.998
$ws.Range("F40").Value = 3.508
$ws.Range("C41").Value = 3.508
$ws.Range("D41").Value = 3.5059999999999998
$ws.Range("E41").Value = 1.9990000000000001
$ws.Range("F41").Value = 1.9930000000000001
$ws.Range("C42").Value = 1.992
$ws.Range("D42").Value = 3.5049999999999999
$ws.Range("E42").Value = 3.508
$ws.Range("F42").Value = 1.9930000000000001
$ws.Range("C43").Value = 1.992
$ws.Range("D43").Value = 1.99
$ws.Range("E43").Value = 3.5110000000000001
$ws.Range("F43").Value = 3.51
$ws.Range("C44").Value = 3.508
$ws.Range("D44").Value = 3.5059999999999998
$ws.Range("E44").Value = 3.51
$ws.Range("F44").Value = 3.5110000000000001
$ws.Range("C45").Value = 1.992
$ws.Range("D45").Value = 3.5070000000000001
$ws.Range("E45").Value = 2
$ws.Range("F45").Value = 3.512

# Sample 5 (rows 47-52)
$ws.Range("C47").Value = 3.4980000000000002
$ws.Range("D47").Value = 1.9890000000000001
$ws.Range("E47").Value = 1.996
$ws.Range("F47").Value = 3.5059999999999998
$ws.Range("C48").Value = 3.4990000000000001
$ws.Range("D48").Value = 3.5019999999999998
$ws.Range("E48").Value = 1.996
$ws.Range("F48").Value = 1.996
$ws.Range("C49").Value = 1.994
$ws.Range("D49").Value = 3.5009999999999999
$ws.Range("E49").Value = 3.5049999999999999
$ws.Range("F49").Value = 1.9950000000000001
$ws.Range("C50").Value = 1.9930000000000001
$ws.Range("D50").Value = 1.9890000000000001
$ws.Range("E50").Value = 3.504
$ws.Range("F50").Value = 3.5070000000000001
$ws.Range("C51").Value = 3.4969999999999999
$ws.Range("D51").Value = 3.5009999999999999
$ws.Range("E51").Value = 3.5030000000000001
$ws.Range("F51").Value = 3.5059999999999998
$ws.Range("C52").Value = 1.994
$ws.Range("D52").Value = 3.5019999999999998
$ws.Range("E52").Value = 1.996
$ws.Range("F52").Value = 3.5070000000000001

# ----------------------------------------------------------------------------
# View state: make "Functionality Test" the active/selected sheet (updates
# workbook.xml activeTab + moves tabSelected from "Execution summary"), and
# set its selected cell to J49.
# ----------------------------------------------------------------------------
$ws.Activate()
$ws.Range("J49").Select()
